# Apply hybrid bold + color highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) across the resume's
# achievements and work-experience bullet points.
#
# Highlight style: bold + font color #2C3E50 (RGB 44,62,80), which as a
# Word wdColor long (BGR-packed) is 0x503E2C = 5258796.

$d = $word.ActiveDocument
$highlightColor = 5258796   # 0x503E2C == RGB(0x2C,0x3E,0x50) in wdColor BGR order

function Format-Metric($para, $term) {
    $r = $para.Range
    $found = $r.Find.Execute($term, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.Font.Bold = 1
        $r.Font.Color = $highlightColor
    } else {
        Write-Host "NOT FOUND: '$term' in paragraph text: $($para.Range.Text)"
    }
    return $found
}

# "• Discovered systematic race coding errors ... accuracy from 23% to 64%"
$p = $d.Paragraphs.Item(10)
Format-Metric $p "23%"
Format-Metric $p "64%"

# "• Utilized advanced sampling methods ... from ±4.2% to ±2.1%, increasing
#   voter turnout prediction accuracy from 71% to 87%, ..."
$p = $d.Paragraphs.Item(12)
Format-Metric $p "±4.2%"
Format-Metric $p "±2.1%"
Format-Metric $p "71%"
Format-Metric $p "87%"

# "• Trigonometric algorithm ... reduced mapping costs by 73.5%, saving
#   campaigns and organizations $4.7M and enabling smaller nonprofits ..."
$p = $d.Paragraphs.Item(13)
Format-Metric $p "73.5%"
Format-Metric $p "`$4.7M"

# "• Built real-time FEC analysis systems ... valued over $2 trillion"
$p = $d.Paragraphs.Item(14)
Format-Metric $p "`$2"

# "• Modernized legacy ETL processes ... reducing processing time by 57%"
$p = $d.Paragraphs.Item(39)
Format-Metric $p "57%"

# "• Algorithmic innovation: Pioneered trigonometric boundary estimation
#   reducing mapping costs 73.5%"
$p = $d.Paragraphs.Item(55)
Format-Metric $p "73.5%"

# "• $4.7M savings enabled nonprofit access"
$p = $d.Paragraphs.Item(56)
Format-Metric $p "`$4.7M"

# "• Platform impact: Built redistricting system serving 12,847 analysts
#   across 89 organizations"
$p = $d.Paragraphs.Item(57)
Format-Metric $p "12,847"
